# Update Mappings 22 Ontologies
# - Fix capitalization: metadata4Ing_* -> metadata4ing_*
# - Add a new "OSMO_DEF" column (F) with default value "[]" for the data row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix header capitalization (B1/C1)
$ws.Range("B1").Value = "metadata4ing_IRI"
$ws.Range("C1").Value = "metadata4ing_DESC"

# New header + data cell for the OSMO_DEF column
$ws.Range("F1").Value = "OSMO_DEF"
$ws.Range("F2").Value = "[]"

# Give the new header cell (F1) the same formatting as the other header cells
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$excel.CutCopyMode = $false
